# Auto-generated Excel COM-interop script
# Applies cached-value corrections to columns H-N across multiple sheets
# as described by the target diff (Sheets/Maduin_Profits.xlsx).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 229.1
$ws.Range("I33").Value = 229.1
$ws.Range("K33").Value = 229.1
$ws.Range("M33").Value = -0.09999999999999432

$ws.Range("H51").Value = 7284
$ws.Range("I51").Value = 6197.6
$ws.Range("K51").Value = 6197.6
$ws.Range("M51").Value = -5713.6

$ws.Range("H106").Value = 1500
$ws.Range("I106").Value = 1500
$ws.Range("K106").Value = 1500
$ws.Range("M106").Value = -869

$ws.Range("H132").Value = 2138.8333
$ws.Range("I132").Value = 2138.8333
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6416.499899999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3886.499899999999
$ws.Range("N132").ClearContents()

$ws.Range("H137").Value = 1215.8334
$ws.Range("I137").Value = 1495.5
$ws.Range("J137").Value = 1159.9
$ws.Range("K137").Value = 4486.5
$ws.Range("L137").Value = 3479.7
$ws.Range("M137").Value = -1936.5
$ws.Range("N137").Value = -8579.700000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 63.8
$ws.Range("I5").Value = 39.666668
$ws.Range("K5").Value = 39.666668
$ws.Range("M5").Value = 72.333332

$ws.Range("H32").Value = 4001.3044
$ws.Range("I32").Value = 3567.1
$ws.Range("J32").Value = 6896
$ws.Range("K32").Value = 3567.1
$ws.Range("L32").Value = 6896
$ws.Range("M32").Value = -3280.1
$ws.Range("N32").Value = -7470

$ws.Range("H37").Value = 7857.143
$ws.Range("J37").Value = 10000
$ws.Range("L37").Value = 10000
$ws.Range("N37").Value = -10546

$ws.Range("H52").Value = 99999
$ws.Range("J52").Value = 99999
$ws.Range("L52").Value = 99999
$ws.Range("N52").Value = -100635

$ws.Range("H61").Value = 2590.111
$ws.Range("I61").Value = 1379.4
$ws.Range("K61").Value = 1379.4
$ws.Range("M61").Value = -1167.4

$ws.Range("H102").Value = 2305.4546
$ws.Range("I102").Value = 2181.5
$ws.Range("J102").Value = 2636
$ws.Range("K102").Value = 2181.5
$ws.Range("L102").Value = 2636
$ws.Range("M102").Value = -559.5
$ws.Range("N102").Value = -5880

$ws.Range("H122").Value = 1951
$ws.Range("I122").Value = 1968.6666
$ws.Range("J122").Value = 1898
$ws.Range("K122").Value = 5905.9998
$ws.Range("L122").Value = 5694
$ws.Range("M122").Value = -3455.9998
$ws.Range("N122").Value = -10594

$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws.Range("H136").Value = 2590.111
$ws.Range("I136").Value = 1379.4
$ws.Range("K136").Value = 4138.200000000001
$ws.Range("M136").Value = -1588.200000000001

$ws.Range("H139").Value = 66443
$ws.Range("J139").Value = 66443
$ws.Range("L139").Value = 66443
$ws.Range("N139").Value = -76723

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 63.8
$ws.Range("I4").Value = 39.666668
$ws.Range("K4").Value = 39.666668
$ws.Range("M4").Value = 75.333332

$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 396.85715
$ws.Range("I7").Value = 106.333336
$ws.Range("J7").Value = 919.8
$ws.Range("K7").Value = 106.333336
$ws.Range("L7").Value = 919.8
$ws.Range("M7").Value = 6.666663999999997
$ws.Range("N7").Value = -1145.8

$ws.Range("H22").Value = 480.7143
$ws.Range("I22").Value = 455
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 455
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -105
$ws.Range("N22").Value = -1200

$ws.Range("H31").Value = 4431.7856
$ws.Range("I31").Value = 3852.125
$ws.Range("J31").Value = 5204.6665
$ws.Range("K31").Value = 3852.125
$ws.Range("L31").Value = 5204.6665
$ws.Range("M31").Value = -3557.125
$ws.Range("N31").Value = -5794.6665

$ws.Range("H34").Value = 4431.7856
$ws.Range("I34").Value = 3852.125
$ws.Range("J34").Value = 5204.6665
$ws.Range("K34").Value = 3852.125
$ws.Range("L34").Value = 5204.6665
$ws.Range("M34").Value = -3650.125
$ws.Range("N34").Value = -5608.6665

$ws.Range("H50").Value = 19833.273
$ws.Range("I50").Value = 11055.333
$ws.Range("J50").Value = 23125
$ws.Range("K50").Value = 11055.333
$ws.Range("L50").Value = 23125
$ws.Range("M50").Value = -10430.333
$ws.Range("N50").Value = -24375

$ws.Range("H51").Value = 21274.75
$ws.Range("J51").Value = 21274.75
$ws.Range("L51").Value = 21274.75
$ws.Range("N51").Value = -22746.75

$ws.Range("H59").Value = 29999.834
$ws.Range("J59").Value = 29999.834
$ws.Range("L59").Value = 29999.834
$ws.Range("N59").Value = -32289.834

$ws.Range("H60").Value = 21666.666
$ws.Range("J60").Value = 21666.666
$ws.Range("L60").Value = 21666.666
$ws.Range("N60").Value = -22688.666

$ws.Range("H61").Value = 21274.75
$ws.Range("J61").Value = 21274.75
$ws.Range("L61").Value = 21274.75
$ws.Range("N61").Value = -21970.75

$ws.Range("H105").Value = 4018
$ws.Range("I105").Value = 2694.5
$ws.Range("J105").Value = 5782.6665
$ws.Range("K105").Value = 2694.5
$ws.Range("L105").Value = 5782.6665
$ws.Range("M105").Value = -947.5
$ws.Range("N105").Value = -9276.666499999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 200000260
$ws.Range("I4").Value = 254.33333
$ws.Range("K4").Value = 762.99999
$ws.Range("M4").Value = -650.99999

$ws.Range("H55").Value = 1632.6666
$ws.Range("I55").Value = 1632.6666
$ws.Range("K55").Value = 4897.9998
$ws.Range("M55").Value = -4720.9998

$ws.Range("H68").Value = 4059
$ws.Range("J68").Value = 4859.6
$ws.Range("L68").Value = 14578.8
$ws.Range("N68").Value = -16200.8

$ws.Range("H71").Value = 4059
$ws.Range("J71").Value = 4859.6
$ws.Range("L71").Value = 43736.4
$ws.Range("N71").Value = -51848.4

$ws.Range("H75").Value = 2400.8572
$ws.Range("I75").Value = 2702
$ws.Range("J75").Value = 1999.3334
$ws.Range("K75").Value = 8106
$ws.Range("L75").Value = 5998.0002
$ws.Range("M75").Value = -7108
$ws.Range("N75").Value = -7994.0002

$ws.Range("H78").Value = 2400.8572
$ws.Range("I78").Value = 2702
$ws.Range("J78").Value = 1999.3334
$ws.Range("K78").Value = 24318
$ws.Range("L78").Value = 17994.0006
$ws.Range("M78").Value = -19326
$ws.Range("N78").Value = -27978.0006

$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H87").Value = 60000
$ws.Range("I87").Value = 60000
$ws.Range("K87").Value = 60000
$ws.Range("M87").Value = -58752

$ws.Range("H90").Value = 60000
$ws.Range("I90").Value = 60000
$ws.Range("K90").Value = 180000
$ws.Range("M90").Value = -173760

$ws.Range("H122").Value = 8800
$ws.Range("I122").Value = 9666.666999999999
$ws.Range("J122").Value = 7500
$ws.Range("K122").Value = 29000.001
$ws.Range("L122").Value = 22500
$ws.Range("M122").Value = -26550.001
$ws.Range("N122").Value = -27400

$ws.Range("H126").Value = 7998.3335
$ws.Range("I126").Value = 7750
$ws.Range("J126").Value = 8495
$ws.Range("K126").Value = 23250
$ws.Range("L126").Value = 25485
$ws.Range("M126").Value = -20780
$ws.Range("N126").Value = -30425

$ws.Range("H132").Value = 4477.3335
$ws.Range("I132").Value = 4338.909
$ws.Range("K132").Value = 13016.727
$ws.Range("M132").Value = -10486.727

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 40000
$ws.Range("J38").Value = 40000
$ws.Range("L38").Value = 40000
$ws.Range("N38").Value = -40820

$ws.Range("H53").Value = 80000
$ws.Range("I53").Value = 80000
$ws.Range("K53").Value = 80000
$ws.Range("M53").Value = -79482

$ws.Range("H61").Value = 1200
$ws.Range("I61").Value = 1200
$ws.Range("K61").Value = 1200
$ws.Range("M61").Value = -998

$ws.Range("H68").Value = 3301
$ws.Range("J68").Value = 4167.6665
$ws.Range("L68").Value = 4167.6665
$ws.Range("N68").Value = -5665.6665

$ws.Range("H71").Value = 3301
$ws.Range("J71").Value = 4167.6665
$ws.Range("L71").Value = 20838.3325
$ws.Range("N71").Value = -28326.3325

$ws.Range("H93").Value = 600
$ws.Range("I93").Value = 600
$ws.Range("K93").Value = 600
$ws.Range("M93").Value = 648

$ws.Range("H113").Value = 1200
$ws.Range("I113").Value = 1200
$ws.Range("K113").Value = 1200
$ws.Range("M113").Value = 970

$ws.Range("H132").Value = 8778.75
$ws.Range("I132").Value = 6999
$ws.Range("J132").Value = 9372
$ws.Range("K132").Value = 20997
$ws.Range("L132").Value = 28116
$ws.Range("M132").Value = -18467
$ws.Range("N132").Value = -33176

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()

$ws.Range("H136").Value = 1539.7142
$ws.Range("I136").Value = 1545.8334
$ws.Range("K136").Value = 4637.5002
$ws.Range("M136").Value = -2087.5002
